$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Ref="D2"; Val="28.381.87"},
    @{Ref="E2"; Val="  -0.39%  "},
    @{Ref="D3"; Val="1.810.51"},
    @{Ref="E3"; Val="  -0.84%  "},
    @{Ref="D4"; Val="1.001"},
    @{Ref="E4"; Val="  -0.22%  "},
    @{Ref="D5"; Val="312.79"},
    @{Ref="E5"; Val="  -1.13%  "},
    @{Ref="D6"; Val="1.000"},
    @{Ref="E6"; Val="  -0.28%  "},
    @{Ref="D7"; Val="0.5159"},
    @{Ref="E7"; Val="  -0.44%  "},
    @{Ref="D8"; Val="0.3990"},
    @{Ref="E8"; Val="  +3.11%  "},
    @{Ref="D9"; Val="0.07878"},
    @{Ref="E9"; Val="  -5.04%  "},
    @{Ref="D10"; Val="1.114"},
    @{Ref="E10"; Val="  -0.91%  "},
    @{Ref="D11"; Val="40.95"},
    @{Ref="E11"; Val="  -2.32%  "},
    @{Ref="D12"; Val="6.375"},
    @{Ref="E12"; Val="  -0.12%  "},
    @{Ref="D13"; Val="1.001"},
    @{Ref="E13"; Val="  -0.20%  "},
    @{Ref="D14"; Val="20.43"},
    @{Ref="E14"; Val="  -3.52%  "},
    @{Ref="D15"; Val="7.335"},
    @{Ref="E15"; Val="  -2.09%  "},
    @{Ref="D16"; Val="1.808.94"},
    @{Ref="E16"; Val="  -1.18%  "},
    @{Ref="D17"; Val="92.78"},
    @{Ref="E17"; Val="  -1.30%  "},
    @{Ref="D18"; Val="0.00001085"},
    @{Ref="E18"; Val="  -3.36%  "},
    @{Ref="D19"; Val="0.06565"},
    @{Ref="E19"; Val="  -1.06%  "},
    @{Ref="D20"; Val="0.9996"},
    @{Ref="E20"; Val="  -0.33%  "},
    @{Ref="D21"; Val="17.33"},
    @{Ref="E21"; Val="  -2.84%  "},
    @{Ref="D22"; Val="6.018"},
    @{Ref="E22"; Val="  -0.84%  "},
    @{Ref="D23"; Val="28.419.75"},
    @{Ref="E23"; Val="  -0.38%  "},
    @{Ref="D24"; Val="11.14"},
    @{Ref="E24"; Val="  -3.01%  "},
    @{Ref="D25"; Val="2.227"},
    @{Ref="E25"; Val="  -0.74%  "},
    @{Ref="D26"; Val="160.90"},
    @{Ref="E26"; Val="  +0.73%  "},
    @{Ref="D27"; Val="20.55"},
    @{Ref="E27"; Val="  -2.73%  "},
    @{Ref="D28"; Val="2.021.58"},
    @{Ref="E28"; Val="  -0.80%  "},
    @{Ref="D29"; Val="2.398"},
    @{Ref="E29"; Val="  -0.67%  "},
    @{Ref="D30"; Val="128.43"},
    @{Ref="E30"; Val="  +1.93%  "},
    @{Ref="D31"; Val="0.1098"},
    @{Ref="E31"; Val="  -0.15%  "},
    @{Ref="D32"; Val="1.069"},
    @{Ref="E32"; Val="  -2.58%  "},
    @{Ref="E33"; Val="  -0.44%  "},
    @{Ref="D34"; Val="5.579"},
    @{Ref="E34"; Val="  -2.80%  "},
    @{Ref="D35"; Val="0.07239"},
    @{Ref="E35"; Val="  -4.69%  "},
    @{Ref="D36"; Val="9.183"},
    @{Ref="E36"; Val="  +4.76%  "},
    @{Ref="D37"; Val="0.02342"},
    @{Ref="E37"; Val="  -1.28%  "},
    @{Ref="D38"; Val="0.2187"},
    @{Ref="E38"; Val="  -2.02%  "},
    @{Ref="D39"; Val="11.63"},
    @{Ref="E39"; Val="  -3.56%  "},
    @{Ref="D40"; Val="5.056"},
    @{Ref="E40"; Val="  -3.90%  "},
    @{Ref="D41"; Val="0.6208"},
    @{Ref="E41"; Val="  -3.22%  "},
    @{Ref="D42"; Val="0.9994"},
    @{Ref="E42"; Val="  -0.34%  "},
    @{Ref="E43"; Val="  -2.75%  "},
    @{Ref="E44"; Val="  -2.72%  "},
    @{Ref="D45"; Val="0.5999"},
    @{Ref="E45"; Val="  -3.49%  "},
    @{Ref="D46"; Val="1.311"},
    @{Ref="E46"; Val="  -6.41%  "},
    @{Ref="D47"; Val="3.732"},
    @{Ref="E47"; Val="  -1.77%  "},
    @{Ref="D48"; Val="125.71"},
    @{Ref="E48"; Val="  -1.81%  "},
    @{Ref="D49"; Val="1.222"},
    @{Ref="E49"; Val="  +1.42%  "},
    @{Ref="D50"; Val="1.931"},
    @{Ref="E50"; Val="  -3.72%  "},
    @{Ref="D51"; Val="0.06841"}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Ref)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Val
    $rng.Style = "Normal"
}
